# The author re-typed the four resistor-header cells so that a space is
# inserted before the trailing "(mV)" unit suffix, e.g. "V_20KΩ(mV)" ->
# "V_20KΩ (mV)". Re-assign the cell text for the four affected headers;
# Excel's shared-string table is rebuilt automatically on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "V_20KΩ (mV)"
$ws.Range("H1").Value = "I_20KΩ (mV)"
$ws.Range("O1").Value = "V_10KΩ (mV)"
$ws.Range("P1").Value = "I_10KΩ (mV)"

# The author's selection ended up on P1 (scrolled right to see the last
# column) instead of the original F5.
$ws.Range("P1").Select()
